$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need to be forced to Text
# so Excel keeps the exact display string (trailing zeros, leading zeros, etc.)
# instead of silently re-parsing it as a numeric value.
$ws.Range("D2").Value = "67.966.38"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "3.749.77"
$ws.Range("E3").Value = "  -2.25%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.66"
$ws.Range("E5").Value = "  -1.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.77"
$ws.Range("E6").Value = "  -3.34%  "
$ws.Range("D7").Value = "3.751.30"
$ws.Range("E7").Value = "  -2.19%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("E10").Value = "  -4.49%  "
$ws.Range("E11").Value = "  -1.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.449"
$ws.Range("E12").Value = "  -2.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000263"
$ws.Range("E13").Value = "  -7.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.99"
$ws.Range("E14").Value = "  -2.50%  "
$ws.Range("D15").Value = "4.376.99"
$ws.Range("E15").Value = "  -2.24%  "
$ws.Range("D16").Value = "3.756.27"
$ws.Range("E16").Value = "  -3.03%  "
$ws.Range("D17").Value = "67.886.78"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.47"
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.05"
$ws.Range("E19").Value = "  -5.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.57"
$ws.Range("E21").Value = "  -2.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "466.11"
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.703"
$ws.Range("E23").Value = "  -3.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.97"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("E25").Value = "  -13.94%  "
$ws.Range("E26").Value = "  -3.41%  "
$ws.Range("E27").Value = "  -1.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.17"
$ws.Range("E28").Value = "  -2.46%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "3.894.29"
$ws.Range("E30").Value = "  -2.37%  "
$ws.Range("E31").Value = "  -2.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.38"
$ws.Range("E32").Value = "  -4.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.90"
$ws.Range("E33").Value = "  -3.57%  "
$ws.Range("E34").Value = "  -4.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.06"
$ws.Range("E35").Value = "  -3.51%  "
$ws.Range("D36").Value = "3.697.15"
$ws.Range("E36").Value = "  -2.64%  "
$ws.Range("E37").Value = "  -2.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.44"
$ws.Range("E38").Value = "  -11.45%  "
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -1.71%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.137"
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.76"
$ws.Range("E41").Value = "  -3.41%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.307"
$ws.Range("E44").Value = "  -3.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.55"
$ws.Range("E45").Value = "  -2.27%  "
$ws.Range("E46").Value = "  -3.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.25"
$ws.Range("E47").Value = "  -2.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "395.72"
$ws.Range("E48").Value = "  -5.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "144.86"
$ws.Range("E49").Value = "  +2.36%  "
$ws.Range("E50").Value = "  -3.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "38.30"
$ws.Range("E51").Value = "  +0.39%  "
